$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 283.53333
$ws.Range("J8").Value = 293.23254
$ws.Range("L8").Value = 879.6976199999999
$ws.Range("N8").Value = -1157.69762

$ws.Range("H17").Value = 1387.2307
$ws.Range("J17").Value = 1387.2307
$ws.Range("L17").Value = 4161.6921
$ws.Range("N17").Value = -4497.6921

$ws.Range("H98").Value = 27029462
$ws.Range("I98").Value = 27780238
$ws.Range("K98").Value = 27780238
$ws.Range("M98").Value = -27778740

$ws.Range("H112").Value = 4550.9395
$ws.Range("J112").Value = 5188.607
$ws.Range("L112").Value = 15565.821
$ws.Range("N112").Value = -17781.821

$ws.Range("H122").Value = 27029462
$ws.Range("I122").Value = 27780238
$ws.Range("K122").Value = 83340714
$ws.Range("M122").Value = -83338264

$ws.Range("H132").Value = 1737.1786
$ws.Range("I132").Value = 1792.7307
$ws.Range("J132").Value = 1015
$ws.Range("K132").Value = 5378.1921
$ws.Range("L132").Value = 3045
$ws.Range("M132").Value = -2848.1921
$ws.Range("N132").Value = -8105

$ws.Range("H135").Value = 834206.7
$ws.Range("I135").Value = 909858.5600000001
$ws.Range("K135").Value = 8188727.040000001
$ws.Range("M135").Value = -8186192.040000001

$ws.Range("H138").Value = 2088572.8
$ws.Range("J138").Value = 2569779.8
$ws.Range("L138").Value = 7709339.399999999
$ws.Range("N138").Value = -7719619.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1317818.8
$ws.Range("I32").Value = 1317818.8
$ws.Range("K32").Value = 1317818.8
$ws.Range("M32").Value = -1317531.8

$ws.Range("H61").Value = 9779.6
$ws.Range("I61").Value = 2204
$ws.Range("K61").Value = 2204
$ws.Range("M61").Value = -1992

$ws.Range("H88").Value = 718.1875
$ws.Range("I88").Value = 367.75
$ws.Range("J88").Value = 1068.625
$ws.Range("K88").Value = 367.75
$ws.Range("L88").Value = 1068.625
$ws.Range("M88").Value = 38.25
$ws.Range("N88").Value = -1880.625

$ws.Range("H91").Value = 718.1875
$ws.Range("I91").Value = 367.75
$ws.Range("J91").Value = 1068.625
$ws.Range("K91").Value = 367.75
$ws.Range("L91").Value = 1068.625
$ws.Range("M91").Value = 1036.25
$ws.Range("N91").Value = -3876.625

$ws.Range("H97").Value = 3789748.5
$ws.Range("I97").Value = 2329.8572
$ws.Range("J97").Value = 10417731
$ws.Range("K97").Value = 2329.8572
$ws.Range("L97").Value = 10417731
$ws.Range("M97").Value = -1833.8572
$ws.Range("N97").Value = -10418723

$ws.Range("H110").Value = 37037628
$ws.Range("J110").Value = 166666820
$ws.Range("L110").Value = 166666820
$ws.Range("N110").Value = -166670910

$ws.Range("H122").Value = 26831.1
$ws.Range("I122").Value = 55077.75
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 165233.25
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -162783.25
$ws.Range("N122").Value = -28900

$ws.Range("H136").Value = 9779.6
$ws.Range("I136").Value = 2204
$ws.Range("K136").Value = 6612
$ws.Range("M136").Value = -4062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11908319
$ws.Range("I20").Value = 20836246
$ws.Range("K20").Value = 20836246
$ws.Range("M20").Value = -20835999

$ws.Range("H86").Value = 35317.8
$ws.Range("I86").Value = 51497.65
$ws.Range("J86").Value = 2958.1
$ws.Range("K86").Value = 51497.65
$ws.Range("L86").Value = 2958.1
$ws.Range("M86").Value = -50374.65
$ws.Range("N86").Value = -5204.1

$ws.Range("H89").Value = 35317.8
$ws.Range("I89").Value = 51497.65
$ws.Range("J89").Value = 2958.1
$ws.Range("K89").Value = 257488.25
$ws.Range("L89").Value = 14790.5
$ws.Range("M89").Value = -251872.25
$ws.Range("N89").Value = -26022.5

$ws.Range("H134").Value = 9622837
$ws.Range("I134").Value = 22729374
$ws.Range("K134").Value = 68188122
$ws.Range("M134").Value = -68185587

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8064.206
$ws.Range("I31").Value = 2894.7144
$ws.Range("J31").Value = 11682.85
$ws.Range("K31").Value = 2894.7144
$ws.Range("L31").Value = 11682.85
$ws.Range("M31").Value = -2599.7144
$ws.Range("N31").Value = -12272.85

$ws.Range("H34").Value = 8064.206
$ws.Range("I34").Value = 2894.7144
$ws.Range("J34").Value = 11682.85
$ws.Range("K34").Value = 2894.7144
$ws.Range("L34").Value = 11682.85
$ws.Range("M34").Value = -2692.7144
$ws.Range("N34").Value = -12086.85

$ws.Range("H58").Value = 8398.678
$ws.Range("I58").Value = 4067.7144
$ws.Range("K58").Value = 4067.7144
$ws.Range("M58").Value = -3864.7144

$ws.Range("H99").Value = 6210.8125
$ws.Range("I99").Value = 6216.091
$ws.Range("J99").Value = 6199.2
$ws.Range("K99").Value = 6216.091
$ws.Range("L99").Value = 6199.2
$ws.Range("M99").Value = -4718.091
$ws.Range("N99").Value = -9195.200000000001

$ws.Range("H126").Value = 6210.8125
$ws.Range("I126").Value = 6216.091
$ws.Range("J126").Value = 6199.2
$ws.Range("K126").Value = 18648.273
$ws.Range("L126").Value = 18597.6
$ws.Range("M126").Value = -16178.273
$ws.Range("N126").Value = -23537.6

$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -40060

$ws.Range("H136").Value = 8398.678
$ws.Range("I136").Value = 4067.7144
$ws.Range("K136").Value = 12203.1432
$ws.Range("M136").Value = -9653.143199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83712.42999999999
$ws.Range("I2").Value = 13090.044
$ws.Range("J2").Value = 315757.44
$ws.Range("K2").Value = 78540.264
$ws.Range("L2").Value = 1894544.64
$ws.Range("M2").Value = -78427.264
$ws.Range("N2").Value = -1894770.64

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H23").Value = 366.41666
$ws.Range("I23").Value = 237.8
$ws.Range("J23").Value = 458.2857
$ws.Range("K23").Value = 713.4000000000001
$ws.Range("L23").Value = 1374.8571
$ws.Range("M23").Value = -478.4000000000001
$ws.Range("N23").Value = -1844.8571

$ws.Range("H34").Value = 5989.2764
$ws.Range("J34").Value = 5923.826
$ws.Range("L34").Value = 17771.478
$ws.Range("N34").Value = -17939.478

$ws.Range("H86").Value = 2003
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 2003
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H116").Value = 2665.6667
$ws.Range("I116").Value = 2247.75
$ws.Range("K116").Value = 6743.25
$ws.Range("M116").Value = -3301.25

$ws.Range("H137").Value = 202192.3
$ws.Range("I137").Value = 145163.86
$ws.Range("K137").Value = 435491.58
$ws.Range("M137").Value = -430391.58

$ws.Range("H138").Value = 80645
$ws.Range("I138").Value = 92580.45
$ws.Range("K138").Value = 277741.35
$ws.Range("M138").Value = -272601.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9128.684999999999
$ws.Range("I70").Value = 5064.857
$ws.Range("K70").Value = 5064.857
$ws.Range("M70").Value = -4794.857

$ws.Range("H73").Value = 9128.684999999999
$ws.Range("I73").Value = 5064.857
$ws.Range("K73").Value = 5064.857
$ws.Range("M73").Value = -4128.857

$ws.Range("H122").Value = 1960054.8
$ws.Range("I122").Value = 2787963.2
$ws.Range("J122").Value = 3180.182
$ws.Range("K122").Value = 8363889.600000001
$ws.Range("L122").Value = 9540.545999999998
$ws.Range("M122").Value = -8361439.600000001
$ws.Range("N122").Value = -14440.546

$ws.Range("H132").Value = 2249.5588
$ws.Range("I132").Value = 2173.9644
$ws.Range("J132").Value = 2602.3333
$ws.Range("K132").Value = 6521.8932
$ws.Range("L132").Value = 7806.999899999999
$ws.Range("M132").Value = -3991.8932
$ws.Range("N132").Value = -12866.9999

$ws.Range("H139").Value = 66660
$ws.Range("J139").Value = 66660
$ws.Range("L139").Value = 66660
$ws.Range("N139").Value = -76940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 199999
$ws.Range("J5").Value = 199999
$ws.Range("L5").Value = 199999
$ws.Range("N5").Value = -200225

$ws.Range("H22").Value = 5451.7393
$ws.Range("I22").Value = 592
$ws.Range("K22").Value = 592
$ws.Range("M22").Value = -297

$ws.Range("H27").Value = 5451.7393
$ws.Range("I27").Value = 592
$ws.Range("K27").Value = 592
$ws.Range("M27").Value = -485

$ws.Range("H46").Value = 4117752.2
$ws.Range("I46").Value = 1788.8
$ws.Range("J46").Value = 6538907
$ws.Range("K46").Value = 1788.8
$ws.Range("L46").Value = 6538907
$ws.Range("M46").Value = -1600.8
$ws.Range("N46").Value = -6539283

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 17548076
$ws.Range("I81").Value = 3467196
$ws.Range("K81").Value = 6934392
$ws.Range("M81").Value = -6933331

$ws.Range("H84").Value = 17548076
$ws.Range("I84").Value = 3467196
$ws.Range("K84").Value = 34671960
$ws.Range("M84").Value = -34666656

$ws.Range("H122").Value = 10959131
$ws.Range("I122").Value = 21915254
$ws.Range("K122").Value = 65745762
$ws.Range("M122").Value = -65743312

$ws.Range("H126").Value = 1681.129
$ws.Range("I126").Value = 1772.238
$ws.Range("J126").Value = 1489.8
$ws.Range("K126").Value = 5316.714
$ws.Range("L126").Value = 4469.4
$ws.Range("M126").Value = -2846.714
$ws.Range("N126").Value = -9409.4

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
